$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text type (so numeric-looking
# strings like "600.24" or "1.00" are preserved as text, not coerced to
# numbers), then restore the cell's style so no stray formatting is left
# behind.
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Column map: A=1 (#), B=2 (Coin), C=3 (Link), D=4 (Price), E=5 (Volume 1h)

# Row 2
Set-TextValue 2 4 "63.931.77"

# Row 3
Set-TextValue 3 4 "3.306.89"
Set-TextValue 3 5 "  +5.93%  "

# Row 4
Set-TextValue 4 5 "  +0.00%  "

# Row 5
Set-TextValue 5 4 "600.24"
Set-TextValue 5 5 "  +0.83%  "

# Row 6
Set-TextValue 6 4 "144.28"
Set-TextValue 6 5 "  +5.63%  "

# Row 8
Set-TextValue 8 4 "3.306.12"
Set-TextValue 8 5 "  +6.08%  "

# Row 9
Set-TextValue 9 5 "  +0.89%  "

# Row 10
Set-TextValue 10 5 "  +3.02%  "

# Row 11
Set-TextValue 11 5 "  +4.94%  "

# Row 12
Set-TextValue 12 5 "  +2.62%  "

# Row 13
Set-TextValue 13 4 "0.0000250"
Set-TextValue 13 5 "  +1.03%  "

# Row 14
Set-TextValue 14 4 "34.94"
Set-TextValue 14 5 "  +2.25%  "

# Row 15
Set-TextValue 15 4 "3.851.05"

# Row 16
Set-TextValue 16 5 "  +1.10%  "

# Row 17
Set-TextValue 17 4 "3.305.17"
Set-TextValue 17 5 "  +5.90%  "

# Row 18
Set-TextValue 18 4 "64.033.59"
Set-TextValue 18 5 "  +1.55%  "

# Row 19
Set-TextValue 19 4 "6.89"
Set-TextValue 19 5 "  +2.34%  "

# Row 20
Set-TextValue 20 4 "481.59"
Set-TextValue 20 5 "  +1.04%  "

# Row 21
Set-TextValue 21 5 "  +0.66%  "

# Row 22
Set-TextValue 22 4 "0.741"
Set-TextValue 22 5 "  +6.15%  "

# Row 23
Set-TextValue 23 4 "8.01"
Set-TextValue 23 5 "  +4.29%  "

# Row 24
Set-TextValue 24 4 "13.57"
Set-TextValue 24 5 "  +4.10%  "

# Row 25
Set-TextValue 25 4 "84.15"
Set-TextValue 25 5 "  -3.57%  "

# Row 27
Set-TextValue 27 4 "2.78"
Set-TextValue 27 5 "  +2.01%  "

# Row 28
Set-TextValue 28 5 "  +1.01%  "

# Row 29 and 30 swap content: RenderToken <-> FirstDigitalUSD
Set-TextValue 29 2 "FirstDigitalUSD"
Set-TextValue 29 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 29 4 "1.00"
Set-TextValue 29 5 "  -0.08%  "

Set-TextValue 30 2 "RenderToken"
Set-TextValue 30 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 30 4 "8.26"
Set-TextValue 30 5 "  +4.36%  "

# Row 31
Set-TextValue 31 4 "2.15"
Set-TextValue 31 5 "  +3.69%  "

# Row 32
Set-TextValue 32 4 "28.45"
Set-TextValue 32 5 "  +5.46%  "

# Row 33
Set-TextValue 33 4 "0.106"
Set-TextValue 33 5 "  -0.65%  "

# Row 34
Set-TextValue 34 5 "  +0.88%  "

# Row 35
Set-TextValue 35 4 "1.11"
Set-TextValue 35 5 "  +2.68%  "

# Row 36
Set-TextValue 36 5 "  +2.79%  "

# Row 37 and 38 swap content: PEPE <-> OKB
Set-TextValue 37 2 "OKB"
Set-TextValue 37 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 37 4 "53.33"
Set-TextValue 37 5 "  +2.58%  "

Set-TextValue 38 2 "PEPE"
Set-TextValue 38 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 38 4 "0.0₃0758"
Set-TextValue 38 5 "  +6.17%  "

# Row 39
Set-TextValue 39 4 "0.0401"
Set-TextValue 39 5 "  +3.29%  "

# Row 40
Set-TextValue 40 4 "431.40"
Set-TextValue 40 5 "  +2.09%  "

# Row 41
Set-TextValue 41 4 "2.81"
Set-TextValue 41 5 "  +4.77%  "

# Row 42
Set-TextValue 42 4 "3.034.55"
Set-TextValue 42 5 "  +5.31%  "

# Row 43
Set-TextValue 43 4 "8.44"
Set-TextValue 43 5 "  +1.98%  "

# Row 44
Set-TextValue 44 5 "  -6.46%  "

# Row 45
Set-TextValue 45 5 "  +0.97%  "

# Row 46
Set-TextValue 46 4 "2.23"
Set-TextValue 46 5 "  +4.31%  "

# Row 47
Set-TextValue 47 4 "26.48"
Set-TextValue 47 5 "  +2.53%  "

# Row 49
Set-TextValue 49 4 "2.33"
Set-TextValue 49 5 "  +1.96%  "

# Row 50
Set-TextValue 50 5 "  +1.67%  "

# Row 51
Set-TextValue 51 4 "35.47"
Set-TextValue 51 5 "  +13.42%  "
